$wb = $excel.ActiveWorkbook

# "Лист3" (Test case 3) is the template for the new "Лист4" (Test case 4) sheet.
$sheet3 = $wb.Worksheets.Item(3)

# Duplicate it, inserting the copy right after Лист3.
$sheet3.Copy([System.Reflection.Missing]::Value, $sheet3)
$sheet4 = $wb.Worksheets.Item(4)
$sheet4.Name = "Лист4"

# Update the test-case specific content on the new sheet.
$sheet4.Range("A2").Value = 4
$sheet4.Range("B2").Value = "Function `ngetHeight in checkedMaze"
$sheet4.Range("C2").Value = "1.Input height"

# Restore the cursor/selection state left on each sheet after editing.
$sheet2 = $wb.Worksheets.Item(2)
$sheet2.Activate()
$sheet2.Range("C20").Select()

$sheet3.Activate()
$sheet3.Range("D28").Select()

$sheet4.Activate()
$sheet4.Range("C24").Select()
